$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-09-19 Friday", $false, $false, $false, $false, $false, $true, 1, $false, "2025-09-20 Saturday", 2) | Out-Null
$d.Content.Find.Execute("18÷9=", $false, $false, $false, $false, $false, $true, 1, $false, "95÷7=", 2) | Out-Null
$d.Content.Find.Execute("55÷6=", $false, $false, $false, $false, $false, $true, 1, $false, "20÷8=", 2) | Out-Null
$d.Content.Find.Execute("55÷2=", $false, $false, $false, $false, $false, $true, 1, $false, "19÷9=", 2) | Out-Null
$d.Content.Find.Execute("18÷4=", $false, $false, $false, $false, $false, $true, 1, $false, "82÷5=", 2) | Out-Null
$d.Content.Find.Execute("26÷3=", $false, $false, $false, $false, $false, $true, 1, $false, "20÷8=", 2) | Out-Null
$d.Content.Find.Execute("94÷4=", $false, $false, $false, $false, $false, $true, 1, $false, "39÷6=", 2) | Out-Null
$d.Content.Find.Execute("24÷5=", $false, $false, $false, $false, $false, $true, 1, $false, "47÷7=", 2) | Out-Null
$d.Content.Find.Execute("50÷9=", $false, $false, $false, $false, $false, $true, 1, $false, "66÷2=", 2) | Out-Null
$d.Content.Find.Execute("58÷9=", $false, $false, $false, $false, $false, $true, 1, $false, "24÷8=", 2) | Out-Null
$d.Content.Find.Execute("22÷8=", $false, $false, $false, $false, $false, $true, 1, $false, "35÷7=", 2) | Out-Null
$d.Content.Find.Execute("10÷7=", $false, $false, $false, $false, $false, $true, 1, $false, "31÷3=", 2) | Out-Null
$d.Content.Find.Execute("32÷2=", $false, $false, $false, $false, $false, $true, 1, $false, "54÷5=", 2) | Out-Null
$d.Content.Find.Execute("70÷7=", $false, $false, $false, $false, $false, $true, 1, $false, "79÷7=", 2) | Out-Null
$d.Content.Find.Execute("40÷4=", $false, $false, $false, $false, $false, $true, 1, $false, "92÷7=", 2) | Out-Null
$d.Content.Find.Execute("50÷2=", $false, $false, $false, $false, $false, $true, 1, $false, "53÷9=", 2) | Out-Null
$d.Content.Find.Execute("98÷9=", $false, $false, $false, $false, $false, $true, 1, $false, "39÷7=", 2) | Out-Null
$d.Content.Find.Execute("21÷7=", $false, $false, $false, $false, $false, $true, 1, $false, "46÷3=", 2) | Out-Null
$d.Content.Find.Execute("71÷7=", $false, $false, $false, $false, $false, $true, 1, $false, "27÷7=", 2) | Out-Null
$d.Content.Find.Execute("94÷8=", $false, $false, $false, $false, $false, $true, 1, $false, "66÷5=", 2) | Out-Null
$d.Content.Find.Execute("66÷8=", $false, $false, $false, $false, $false, $true, 1, $false, "41÷7=", 2) | Out-Null
$d.Content.Find.Execute("95÷5=", $false, $false, $false, $false, $false, $true, 1, $false, "45÷4=", 2) | Out-Null
$d.Content.Find.Execute("40÷6=", $false, $false, $false, $false, $false, $true, 1, $false, "75÷4=", 2) | Out-Null
$d.Content.Find.Execute("55÷5=", $false, $false, $false, $false, $false, $true, 1, $false, "55÷6=", 2) | Out-Null
$d.Content.Find.Execute("94÷9=", $false, $false, $false, $false, $false, $true, 1, $false, "94÷5=", 2) | Out-Null
$d.Content.Find.Execute("48÷9=", $false, $false, $false, $false, $false, $true, 1, $false, "14÷8=", 2) | Out-Null
